# Delete unused variables in Molten Salt Tower Parasitics UI page
# Adds rows describing deleted variables to the "SAM Variable Changes" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SAM Variable Changes")

# Row 36-37 reuse the formatting of row 35 (F column keeps the "25" style).
$ws.Range("A35:H35").Copy()
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("A37").PasteSpecial(-4122)

# Rows 38-42 reuse the formatting of row 34 (F column uses the "19" style).
$ws.Range("A34:H34").Copy()
$ws.Range("A38").PasteSpecial(-4122)
$ws.Range("A39").PasteSpecial(-4122)
$ws.Range("A40").PasteSpecial(-4122)
$ws.Range("A41").PasteSpecial(-4122)
$ws.Range("A42").PasteSpecial(-4122)

# Column D is blank in both template rows; drop the stray empty cell that
# the whole-row paste leaves behind so rows 36-42 stay D-less like the rest.
$ws.Range("D36:D42").ClearContents()

$rows = @(
    @{ Row = 36; C = "P_storage_pump";    F = "storage HTF = rec/pc HTF (no storage HX), so no pumping losses" },
    @{ Row = 37; C = "storage_bypass";    F = "storage HTF = rec/pc HTF (no storage HX), so no pumping losses" },
    @{ Row = 38; C = "recirc_source";     F = "not used" },
    @{ Row = 39; C = "recirc_htf_eff";    F = "not used" },
    @{ Row = 40; C = "flow_from_storage"; F = "not used" },
    @{ Row = 41; C = "P_hot_tank";        F = "not used" },
    @{ Row = 42; C = "csp.pt.par.bop_c1"; F = "not used" }
)

foreach ($entry in $rows) {
    $r = $entry.Row

    $ws.Cells.Item($r, 1).Value = "Deleted variable"
    $ws.Cells.Item($r, 2).Value = "number"
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 5).Value = "Molten Salt Tower Parasitics"
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = "N/A"
    $ws.Cells.Item($r, 8).Value = "Ty"
}

$ws.Application.CutCopyMode = $false

# Scroll the view down a bit and leave the selection on the row right after
# the newly added data, matching where the editor ended up after the edit.
$ws.Activate()
$ws.Range("A43").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
